$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Mon Jul 17 21:11:47 EDT 2023"
$ws.Range("B3").Value = "Mon Jul 17 21:11:57 EDT 2023"
$ws.Range("B4").Value = "Mon Jul 17 21:12:07 EDT 2023"
